$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "Anno di Produzione"
$ws.Range("E2").Select()
